$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.891.69"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.219.39"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "2.562.63"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "2.231.03"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.732"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "39.824.39"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.47%  "
$ws.Range("D42").Value = "2.104.01"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0272"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("E47").Value = "  -7.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("D49").Value = "2.437.36"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.58%  "
